{"js": "// Locate the existing bullet paragraph that ends the \"V1.1 (unreleased)\" list:\n// \"All functions now work due to poor pixel value correctors (ie. eliminated\n// NaN and infinity values from images)\". Two new bullet items are added right\n// after it (and therefore right before the \"Fluidics\" heading), using the\n// same list (numId 13) / style (List Paragraph) as their neighbours:\n//   - \"Optional rolling ball background subtraction\" + \" (default)\"\n//   - \"Linear blend in focus tile pieces\"\n\nconst anchorText = \"All functions now work due to poor pixel value correctors\";\nconst searchResults = context.document.body.search(anchorText, { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found: \" + anchorText);\n}\n\n// Resolve the whole paragraph that contains the matched text.\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\nanchorParagraph.load(\"text\");\nawait context.sync();\n\n// Insert the first new bullet right after the anchor paragraph. insertParagraph\n// copies the anchor's paragraph properties (style \"List Paragraph\" + the\n// numId 13 numbering), matching the target markup.\nconst ballParagraph = anchorParagraph.insertParagraph(\n  \"Optional rolling ball background subtraction\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Append the \" (default)\" suffix onto that same bullet.\nballParagraph.insertText(\" (default)\", Word.InsertLocation.end);\nawait context.sync();\n\n// Insert the second new bullet right after the first new one.\nballParagraph.insertParagraph(\n  \"Linear blend in focus tile pieces\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Locate the existing bullet paragraph that ends the \"V1.1 (unreleased)\" list:\n# \"All functions now work due to poor pixel value correctors (ie. eliminated\n# NaN and infinity values from images)\". Two new bullet items are added right\n# after it (and therefore right before the \"Fluidics\" heading), using the\n# same list (numId 13) / style (List Paragraph) as their neighbours:\n#   - \"Optional rolling ball background subtraction\" + \" (default)\"\n#   - \"Linear blend in focus tile pieces\"\n\n$d = $word.ActiveDocument\n\n$anchor = $d.Content\n$find = $anchor.Find\n$find.ClearFormatting()\n$find.Text = \"All functions now work due to poor pixel value correctors\"\n$find.Forward = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Expand the found hit to the full paragraph (style/numbering live on the\n# paragraph mark, so InsertParagraphAfter below inherits \"List Paragraph\" /\n# numId 13 from this paragraph automatically).\n$anchor.Expand(4) | Out-Null  # wdParagraph = 4\n$anchor.InsertParagraphAfter()\n\n# The newly-minted (empty) paragraph sits right after $anchor's end; setting\n# .Text on the collapsed range there fills it in without touching the\n# paragraph mark / inherited numPr.\n$ballRange = $d.Range($anchor.End, $anchor.End)\n$ballRange.Text = \"Optional rolling ball background subtraction\"\n\n# Append the \" (default)\" suffix onto that same bullet.\n$suffixRange = $d.Range($ballRange.End, $ballRange.End)\n$suffixRange.Text = \" (default)\"\n\n# Insert the second new bullet right after the first new one.\n$ballParagraph = $d.Range($suffixRange.End, $suffixRange.End)\n$ballParagraph.Expand(4) | Out-Null  # wdParagraph = 4\n$ballParagraph.InsertParagraphAfter()\n\n$blendRange = $d.Range($ballParagraph.End, $ballParagraph.End)\n$blendRange.Text = \"Linear blend in focus tile pieces\"\n"}
